# Feature Iteration and Bug Fixing
#
# 1) Fix the JSON schema sample string stored in B2 ("int" -> "integer").
# 2) Move/collapse the saved sheet selection to just B2 (was A4:D4 D9 / D9).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("namespace")

$ws.Range("B2").Value = '[{"name":"t0","type":"integer"},{"name":"t1","type":"long"},{"name":"t2","type":"float"}]'

$ws.Range("B2").Select()
